$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: Update the "总计" (grand-total) sheet.
#   Insert a "2022-Q4" row above the existing "2020-Q4" row, so the existing
#   row moves from row 2 down to row 3 and the new row becomes row 2.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Read the existing "2020-Q4" row values (use Value2 - Value's getter is not
# reliable for reads in this host, only as a setter).
$oldB = $wsTotal.Cells.Item(2, 2).Value2
$oldC = $wsTotal.Cells.Item(2, 3).Value2
$oldD = $wsTotal.Cells.Item(2, 4).Value2

# Duplicate row 2's formatting onto row 3 first, so the relocated row keeps
# the same look (bold/centered/bordered index style, etc.). -4122 = xlPasteFormats.
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)

# Put the original "2020-Q4" values back, now on row 3. Column A is a
# 0-based row-position index, so it becomes 1 (second data row) here.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = $oldB
$wsTotal.Cells.Item(3, 3).Value = $oldC
$wsTotal.Cells.Item(3, 4).Value = $oldD

# Write the new "2022-Q4" totals into row 2.
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.03

# ---------------------------------------------------------------------------
# Step 2: Insert a new "2022-Q4" holdings sheet between "总计" and "2020-Q4".
#   Activating "2020-Q4" first makes Worksheets.Add() (with no explicit
#   Before/After) drop the new sheet immediately in front of it.
# ---------------------------------------------------------------------------
$ws2020 = $wb.Worksheets.Item("2020-Q4")
$ws2020.Activate()

$ws2022 = $wb.Worksheets.Add()
$ws2022.Name = "2022-Q4"

# Header row.
$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

# Force text formatting for columns that hold numeric-looking text (fund
# code, fund size, position %, etc.) so leading zeros / exact text is kept
# instead of Excel silently converting them to numbers. The format is reset
# back to "Normal" afterwards so the cells end up with the plain/default
# style again (only their underlying value type stays text).
$ws2022.Range("B2:B3").NumberFormat = "@"
$ws2022.Range("D2:G3").NumberFormat = "@"

# Row 2: 014133 / 工银中证500六个月持有指数增强A
$ws2022.Cells.Item(2, 1).Value = 0
$ws2022.Cells.Item(2, 2).Value = "014133"
$ws2022.Cells.Item(2, 3).Value = "工银中证500六个月持有指数增强A"
$ws2022.Cells.Item(2, 4).Value = "1.59"
$ws2022.Cells.Item(2, 5).Value = "94.29"
$ws2022.Cells.Item(2, 6).Value = "1.31"
$ws2022.Cells.Item(2, 7).Value = "0.0208"
$ws2022.Cells.Item(2, 8).Value = 6

# Row 3: 014134 / 工银中证500六个月持有指数增强C
$ws2022.Cells.Item(3, 1).Value = 1
$ws2022.Cells.Item(3, 2).Value = "014134"
$ws2022.Cells.Item(3, 3).Value = "工银中证500六个月持有指数增强C"
$ws2022.Cells.Item(3, 4).Value = "0.88"
$ws2022.Cells.Item(3, 5).Value = "94.29"
$ws2022.Cells.Item(3, 6).Value = "1.31"
$ws2022.Cells.Item(3, 7).Value = "0.0115"
$ws2022.Cells.Item(3, 8).Value = 6

# Drop the temporary "@" text format now that the values are locked in as
# text, restoring the plain/default cell style used by the source data.
$ws2022.Range("B2:B3").Style = "Normal"
$ws2022.Range("D2:G3").Style = "Normal"

# ---------------------------------------------------------------------------
# Step 3: Match the look & feel of the other sheets (bold/centered/bordered
# header row and "index" column), by copying formatting from "总计".
# ---------------------------------------------------------------------------
# -4122 = xlPasteFormats, so only styles move across (values are untouched).
$wsTotal.Cells.Item(1, 2).Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Cells.Item(2, 1).Copy()
$ws2022.Cells.Item(2, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(2, 1).Copy()
$ws2022.Cells.Item(3, 1).PasteSpecial(-4122)

$ws2022.Activate()
$ws2022.Range("A1").Select() | Out-Null

Write-Host "Done: added 2022-Q4 sheet and updated summary totals"
